# Re-process metadata sheet with newly curated dimensions:
# Columns C (residencia-comarca-nombre), E (edad-grupos-quinquenales),
# F (residencia-ccaa-nombre), K (sexo) and M (residencia-provincia-nombre)
# move from "iaest-dimension:*" to "iaest-measure:*" (row 2), from "dim" to
# "medida" (row 3) and from their bespoke type (skos:Concept / URI-*) to
# "xsd:int" (row 4). Column B (ano) remains the only true dimension.
# The now-unused per-dimension mapping workbook references in row 5
# (edad-grupos-quinquenales, residencia-ccaa-nombre, sexo) are removed;
# only the "ano" mapping (column B) remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: sdmx/iaest annotation -> iaest-measure:<field>
$ws.Range("C2").Value = "iaest-measure:residencia-comarca-nombre"
$ws.Range("E2").Value = "iaest-measure:edad-grupos-quinquenales"
$ws.Range("F2").Value = "iaest-measure:residencia-ccaa-nombre"
$ws.Range("K2").Value = "iaest-measure:sexo"
$ws.Range("M2").Value = "iaest-measure:residencia-provincia-nombre"

# Row 3: dim -> medida for every curated measure column
$ws.Range("C3").Value = "medida"
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "medida"
$ws.Range("K3").Value = "medida"
$ws.Range("M3").Value = "medida"

# Row 4: type column -> xsd:int for every curated measure column
$ws.Range("C4").Value = "xsd:int"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("K4").Value = "xsd:int"
$ws.Range("M4").Value = "xsd:int"

# Row 5: drop the mapping-file references that no longer apply
$ws.Range("E5").Clear()
$ws.Range("F5").Clear()
$ws.Range("K5").Clear()
